# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Kazajistan/Egipto order (rows 29 & 30 in column A) ---
$ws.Range("A29").Value2 = "Kazajistan"
$ws.Range("A30").Value2 = "Egipto"

# --- Update "Datos actualizados" timestamp in A1 ---
$ws.Range("A1").Value2 = "Datos actualizados a 5 de Agosto de 2020 a las 05:31"

# --- Update numeric data for affected country rows ---

# Row 29 (Kazajistan)
$ws.Range("B29").Value2 = 94882
$ws.Range("C29").Value2 = 1062
$ws.Range("D29").Value2 = 67031
$ws.Range("E29").Value2 = 26793
$ws.Range("H29").Value2 = 1058

# Row 30 (Egipto)
$ws.Range("B30").Value2 = 94752
$ws.Range("C30").Value2 = 0
$ws.Range("D30").Value2 = 45569
$ws.Range("E30").Value2 = 44271
$ws.Range("H30").Value2 = 4912

# Row 39
$ws.Range("B39").Value2 = 70648
$ws.Range("C39").Value2 = 334
$ws.Range("D39").Value2 = 17639
$ws.Range("E39").Value2 = 43157
$ws.Range("G39").Value2 = 2
$ws.Range("H39").Value2 = 9852

# Row 51
$ws.Range("B51").Value2 = 44299
$ws.Range("C51").Value2 = 505
$ws.Range("D51").Value2 = 5921
$ws.Range("E51").Value2 = 36978
$ws.Range("G51").Value2 = 16
$ws.Range("H51").Value2 = 1400

# Row 72
$ws.Range("B72").Value2 = 19444
$ws.Range("C72").Value2 = 714
$ws.Range("D72").Value2 = 10799
$ws.Range("E72").Value2 = 8398
$ws.Range("G72").Value2 = 15
$ws.Range("H72").Value2 = 247

# Row 112
$ws.Range("E112").Value2 = 1487
$ws.Range("G112").Value2 = 1
$ws.Range("H112").Value2 = 42

# Row 176
$ws.Range("B176").Value2 = 243
$ws.Range("C176").Value2 = 3
$ws.Range("D176").Value2 = 202
$ws.Range("E176").Value2 = 41
